# edit.ps1
# Scheduled market-data refresh for the Yojimbo Profits workbook.
# Updates cached item-price / leve-profit figures (columns H-N) on a
# handful of rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR
# sheets to reflect the latest market board snapshot. Column layout:
#   H currentAveragePrice      I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ              L LevePriceHQ
#   M LeveProfitNQ             N LeveProfitHQ

# Auto-generated cell updates grouped by sheet/row
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: H32=686.6316, I32=528.5, J32=957.7143, K32=528.5, L32=957.7143, M32=-202.5, N32=-1609.7143
$ws.Range("H32").Value = 686.6316
$ws.Range("I32").Value = 528.5
$ws.Range("J32").Value = 957.7143
$ws.Range("K32").Value = 528.5
$ws.Range("L32").Value = 957.7143
$ws.Range("M32").Value = -202.5
$ws.Range("N32").Value = -1609.7143

# Row 55: H55=89.545456, I55=60, K55=60, M55=154
$ws.Range("H55").Value = 89.545456
$ws.Range("I55").Value = 60
$ws.Range("K55").Value = 60
$ws.Range("M55").Value = 154

# Row 116: H116=13893000, J116=41671584, L116=41671584, N116=-41678468
$ws.Range("H116").Value = 13893000
$ws.Range("J116").Value = 41671584
$ws.Range("L116").Value = 41671584
$ws.Range("N116").Value = -41678468

# Row 125: H125=2391.5557, I125=2210.6667, J125=2482, K125=19896.0003, L125=22338, M125=-17436.0003, N125=-27258
$ws.Range("H125").Value = 2391.5557
$ws.Range("I125").Value = 2210.6667
$ws.Range("J125").Value = 2482
$ws.Range("K125").Value = 19896.0003
$ws.Range("L125").Value = 22338
$ws.Range("M125").Value = -17436.0003
$ws.Range("N125").Value = -27258

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32=2665.67, I32=2414.699, K32=2414.699, M32=-2127.699
$ws.Range("H32").Value = 2665.67
$ws.Range("I32").Value = 2414.699
$ws.Range("K32").Value = 2414.699
$ws.Range("M32").Value = -2127.699

# Row 74: H74=2227.348, I74=1516.6666, J74=2684.2144, K74=1516.6666, L74=2684.2144, M74=-642.6666, N74=-4432.2144
$ws.Range("H74").Value = 2227.348
$ws.Range("I74").Value = 1516.6666
$ws.Range("J74").Value = 2684.2144
$ws.Range("K74").Value = 1516.6666
$ws.Range("L74").Value = 2684.2144
$ws.Range("M74").Value = -642.6666
$ws.Range("N74").Value = -4432.2144

# Row 77: H77=2227.348, I77=1516.6666, J77=2684.2144, K77=7583.333000000001, L77=13421.072, M77=-3215.333000000001, N77=-22157.072
$ws.Range("H77").Value = 2227.348
$ws.Range("I77").Value = 1516.6666
$ws.Range("J77").Value = 2684.2144
$ws.Range("K77").Value = 7583.333000000001
$ws.Range("L77").Value = 13421.072
$ws.Range("M77").Value = -3215.333000000001
$ws.Range("N77").Value = -22157.072

# Row 122: H122=2281.7036, I122=2178.9412, J122=2456.4, K122=6536.823600000001, L122=7369.200000000001, M122=-4086.823600000001, N122=-12269.2
$ws.Range("H122").Value = 2281.7036
$ws.Range("I122").Value = 2178.9412
$ws.Range("J122").Value = 2456.4
$ws.Range("K122").Value = 6536.823600000001
$ws.Range("L122").Value = 7369.200000000001
$ws.Range("M122").Value = -4086.823600000001
$ws.Range("N122").Value = -12269.2

# Row 132: H132=2964.5532, I132=2743.0857, J132=3610.5, K132=8229.257100000001, L132=10831.5, M132=-5699.257100000001, N132=-15891.5
$ws.Range("H132").Value = 2964.5532
$ws.Range("I132").Value = 2743.0857
$ws.Range("J132").Value = 3610.5
$ws.Range("K132").Value = 8229.257100000001
$ws.Range("L132").Value = 10831.5
$ws.Range("M132").Value = -5699.257100000001
$ws.Range("N132").Value = -15891.5

$ws = $wb.Worksheets.Item("BSM")
# Row 64: H64=567.6, I64=648, J64=555.2308, K64=648, L64=555.2308, M64=-423, N64=-1005.2308
$ws.Range("H64").Value = 567.6
$ws.Range("I64").Value = 648
$ws.Range("J64").Value = 555.2308
$ws.Range("K64").Value = 648
$ws.Range("L64").Value = 555.2308
$ws.Range("M64").Value = -423
$ws.Range("N64").Value = -1005.2308

# Row 67: H67=567.6, I67=648, J67=555.2308, K67=648, L67=555.2308, M67=132, N67=-2115.2308
$ws.Range("H67").Value = 567.6
$ws.Range("I67").Value = 648
$ws.Range("J67").Value = 555.2308
$ws.Range("K67").Value = 648
$ws.Range("L67").Value = 555.2308
$ws.Range("M67").Value = 132
$ws.Range("N67").Value = -2115.2308

# Row 94: H94=1073.5454, I94=930.9, J94=2500, K94=930.9, L94=2500, M94=-479.9, N94=-3402
$ws.Range("H94").Value = 1073.5454
$ws.Range("I94").Value = 930.9
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 930.9
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -479.9
$ws.Range("N94").Value = -3402

# Row 105: H105=1851.2667, I105=1772.4166, K105=1772.4166, M105=-25.41660000000002
$ws.Range("H105").Value = 1851.2667
$ws.Range("I105").Value = 1772.4166
$ws.Range("K105").Value = 1772.4166
$ws.Range("M105").Value = -25.41660000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 58: H58=1887.1143, I58=2182.913, J58=1320.1666, K58=2182.913, L58=1320.1666, M58=-1979.913, N58=-1726.1666
$ws.Range("H58").Value = 1887.1143
$ws.Range("I58").Value = 2182.913
$ws.Range("J58").Value = 1320.1666
$ws.Range("K58").Value = 2182.913
$ws.Range("L58").Value = 1320.1666
$ws.Range("M58").Value = -1979.913
$ws.Range("N58").Value = -1726.1666

# Row 105: H105=1900.125, I105=1900.125, K105=1900.125, M105=-153.125
$ws.Range("H105").Value = 1900.125
$ws.Range("I105").Value = 1900.125
$ws.Range("K105").Value = 1900.125
$ws.Range("M105").Value = -153.125

# Row 136: H136=1887.1143, I136=2182.913, J136=1320.1666, K136=6548.739, L136=3960.4998, M136=-3998.739, N136=-9060.4998
$ws.Range("H136").Value = 1887.1143
$ws.Range("I136").Value = 2182.913
$ws.Range("J136").Value = 1320.1666
$ws.Range("K136").Value = 6548.739
$ws.Range("L136").Value = 3960.4998
$ws.Range("M136").Value = -3998.739
$ws.Range("N136").Value = -9060.4998

$ws = $wb.Worksheets.Item("CUL")
# Row 132: H132=1371.6666, I132=965.6667, J132=1534.0667, K132=8691.0003, L132=13806.6003, M132=-6161.0003, N132=-18866.6003
$ws.Range("H132").Value = 1371.6666
$ws.Range("I132").Value = 965.6667
$ws.Range("J132").Value = 1534.0667
$ws.Range("K132").Value = 8691.0003
$ws.Range("L132").Value = 13806.6003
$ws.Range("M132").Value = -6161.0003
$ws.Range("N132").Value = -18866.6003

$ws = $wb.Worksheets.Item("GSM")
# Row 122: H122=2102.4194, I122=1407.75, J122=3365.4546, K122=4223.25, L122=10096.3638, M122=-1773.25, N122=-14996.3638
$ws.Range("H122").Value = 2102.4194
$ws.Range("I122").Value = 1407.75
$ws.Range("J122").Value = 3365.4546
$ws.Range("K122").Value = 4223.25
$ws.Range("L122").Value = 10096.3638
$ws.Range("M122").Value = -1773.25
$ws.Range("N122").Value = -14996.3638

$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22=858.625, I22=330, J22=934.1429000000001, K22=330, L22=934.1429000000001, M22=-35, N22=-1524.1429
$ws.Range("H22").Value = 858.625
$ws.Range("I22").Value = 330
$ws.Range("J22").Value = 934.1429000000001
$ws.Range("K22").Value = 330
$ws.Range("L22").Value = 934.1429000000001
$ws.Range("M22").Value = -35
$ws.Range("N22").Value = -1524.1429

# Row 27: H27=858.625, I27=330, J27=934.1429000000001, K27=330, L27=934.1429000000001, M27=-223, N27=-1148.1429
$ws.Range("H27").Value = 858.625
$ws.Range("I27").Value = 330
$ws.Range("J27").Value = 934.1429000000001
$ws.Range("K27").Value = 330
$ws.Range("L27").Value = 934.1429000000001
$ws.Range("M27").Value = -223
$ws.Range("N27").Value = -1148.1429

# Row 40: H40=2724, I40=2529.5293, J40=3275, K40=2529.5293, L40=3275, M40=-2393.5293, N40=-3547
$ws.Range("H40").Value = 2724
$ws.Range("I40").Value = 2529.5293
$ws.Range("J40").Value = 3275
$ws.Range("K40").Value = 2529.5293
$ws.Range("L40").Value = 3275
$ws.Range("M40").Value = -2393.5293
$ws.Range("N40").Value = -3547

# Row 46: H46=207397.42, I46=2071.4285, J46=327170.9, K46=2071.4285, L46=327170.9, M46=-1883.4285, N46=-327546.9
$ws.Range("H46").Value = 207397.42
$ws.Range("I46").Value = 2071.4285
$ws.Range("J46").Value = 327170.9
$ws.Range("K46").Value = 2071.4285
$ws.Range("L46").Value = 327170.9
$ws.Range("M46").Value = -1883.4285
$ws.Range("N46").Value = -327546.9

# Row 61: H61=2574.75, I61=2166.3333, J61=3800, K61=2166.3333, L61=3800, M61=-1964.3333, N61=-4204
$ws.Range("H61").Value = 2574.75
$ws.Range("I61").Value = 2166.3333
$ws.Range("J61").Value = 3800
$ws.Range("K61").Value = 2166.3333
$ws.Range("L61").Value = 3800
$ws.Range("M61").Value = -1964.3333
$ws.Range("N61").Value = -4204

# Row 113: H113=2574.75, I113=2166.3333, J113=3800, K113=2166.3333, L113=3800, M113=3.666700000000219, N113=-8140
$ws.Range("H113").Value = 2574.75
$ws.Range("I113").Value = 2166.3333
$ws.Range("J113").Value = 3800
$ws.Range("K113").Value = 2166.3333
$ws.Range("L113").Value = 3800
$ws.Range("M113").Value = 3.666700000000219
$ws.Range("N113").Value = -8140

# Row 122: H122=3962.2666, I122=4209.25, J122=3680, K122=12627.75, L122=11040, M122=-10177.75, N122=-15940
$ws.Range("H122").Value = 3962.2666
$ws.Range("I122").Value = 4209.25
$ws.Range("J122").Value = 3680
$ws.Range("K122").Value = 12627.75
$ws.Range("L122").Value = 11040
$ws.Range("M122").Value = -10177.75
$ws.Range("N122").Value = -15940

$ws = $wb.Worksheets.Item("WVR")
# Row 107: H107=274.53333, I107=251.33333, K107=753.99999, M107=1166.00001
$ws.Range("H107").Value = 274.53333
$ws.Range("I107").Value = 251.33333
$ws.Range("K107").Value = 753.99999
$ws.Range("M107").Value = 1166.00001

# Row 113: H113=264, I113=252.875, J113=353, K113=758.625, L113=1059, M113=1411.375, N113=-5399
$ws.Range("H113").Value = 264
$ws.Range("I113").Value = 252.875
$ws.Range("J113").Value = 353
$ws.Range("K113").Value = 758.625
$ws.Range("L113").Value = 1059
$ws.Range("M113").Value = 1411.375
$ws.Range("N113").Value = -5399

# Row 122: H122=358354.2, I122=667292.1, J122=1887.3077, K122=2001876.3, L122=5661.9231, M122=-1999426.3, N122=-10561.9231
$ws.Range("H122").Value = 358354.2
$ws.Range("I122").Value = 667292.1
$ws.Range("J122").Value = 1887.3077
$ws.Range("K122").Value = 2001876.3
$ws.Range("L122").Value = 5661.9231
$ws.Range("M122").Value = -1999426.3
$ws.Range("N122").Value = -10561.9231

# Row 126: H126=334392.28, I126=714875.9399999999, J126=1469.0625, K126=2144627.82, L126=4407.1875, M126=-2142157.82, N126=-9347.1875
$ws.Range("H126").Value = 334392.28
$ws.Range("I126").Value = 714875.9399999999
$ws.Range("J126").Value = 1469.0625
$ws.Range("K126").Value = 2144627.82
$ws.Range("L126").Value = 4407.1875
$ws.Range("M126").Value = -2142157.82
$ws.Range("N126").Value = -9347.1875

# Row 132: H132=831.0658, I132=610.9643, J132=1447.35, K132=1832.8929, L132=4342.049999999999, M132=697.1071000000002, N132=-9402.049999999999
$ws.Range("H132").Value = 831.0658
$ws.Range("I132").Value = 610.9643
$ws.Range("J132").Value = 1447.35
$ws.Range("K132").Value = 1832.8929
$ws.Range("L132").Value = 4342.049999999999
$ws.Range("M132").Value = 697.1071000000002
$ws.Range("N132").Value = -9402.049999999999

# Row 136: H136=738.40814, I136=569.79486, J136=1396, K136=1709.38458, L136=4188, M136=840.6154200000001, N136=-9288
$ws.Range("H136").Value = 738.40814
$ws.Range("I136").Value = 569.79486
$ws.Range("J136").Value = 1396
$ws.Range("K136").Value = 1709.38458
$ws.Range("L136").Value = 4188
$ws.Range("M136").Value = 840.6154200000001
$ws.Range("N136").Value = -9288

